# Auto-generated edit script: updates Leve profit-calculation values
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR crafting-class sheets
# (scheduled-runner price refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4973.3125
$ws.Range("I98").Value = 2707.84
$ws.Range("J98").Value = 7435.7827
$ws.Range("K98").Value = 2707.84
$ws.Range("L98").Value = 7435.7827
$ws.Range("M98").Value = -1209.84
$ws.Range("N98").Value = -10431.7827
$ws.Range("H122").Value = 4973.3125
$ws.Range("I122").Value = 2707.84
$ws.Range("J122").Value = 7435.7827
$ws.Range("K122").Value = 8123.52
$ws.Range("L122").Value = 22307.3481
$ws.Range("M122").Value = -5673.52
$ws.Range("N122").Value = -27207.3481
$ws.Range("H127").Value = 1391.5862
$ws.Range("J127").Value = 1595.238
$ws.Range("L127").Value = 4785.714
$ws.Range("N127").Value = -14705.714
$ws.Range("H138").Value = 2989.8113
$ws.Range("I138").Value = 1518
$ws.Range("J138").Value = 3881.818
$ws.Range("K138").Value = 4554
$ws.Range("L138").Value = 11645.454
$ws.Range("M138").Value = 586
$ws.Range("N138").Value = -21925.454
$ws.Range("H141").Value = 12125.546
$ws.Range("I141").Value = 17153.715
$ws.Range("K141").Value = 51461.145
$ws.Range("M141").Value = -46281.145

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4495.0967
$ws.Range("I32").Value = 4471.8086
$ws.Range("K32").Value = 4471.8086
$ws.Range("M32").Value = -4184.8086
$ws.Range("H37").Value = 30579.46
$ws.Range("J37").Value = 30091.625
$ws.Range("L37").Value = 30091.625
$ws.Range("N37").Value = -30637.625
$ws.Range("H44").Value = 49571.43
$ws.Range("J44").Value = 49571.43
$ws.Range("L44").Value = 49571.43
$ws.Range("N44").Value = -50547.43
$ws.Range("H45").Value = 2193
$ws.Range("I45").Value = 2407.7144
$ws.Range("K45").Value = 2407.7144
$ws.Range("M45").Value = -2030.7144
$ws.Range("H55").Value = 50000
$ws.Range("J55").Value = 50000
$ws.Range("L55").Value = 50000
$ws.Range("N55").Value = -50630
$ws.Range("H63").Value = 9898293
$ws.Range("I63").Value = 27704460
$ws.Range("J63").Value = 5977.778
$ws.Range("K63").Value = 27704460
$ws.Range("L63").Value = 5977.778
$ws.Range("M63").Value = -27703774
$ws.Range("N63").Value = -7349.778
$ws.Range("H66").Value = 9898293
$ws.Range("I66").Value = 27704460
$ws.Range("J66").Value = 5977.778
$ws.Range("K66").Value = 138522300
$ws.Range("L66").Value = 29888.89
$ws.Range("M66").Value = -138518868
$ws.Range("N66").Value = -36752.89
$ws.Range("H74").Value = 297393.6
$ws.Range("I74").Value = 485117.5
$ws.Range("K74").Value = 485117.5
$ws.Range("M74").Value = -484243.5
$ws.Range("H77").Value = 297393.6
$ws.Range("I77").Value = 485117.5
$ws.Range("K77").Value = 2425587.5
$ws.Range("M77").Value = -2421219.5
$ws.Range("H80").Value = 34772.8
$ws.Range("J80").Value = 34772.8
$ws.Range("L80").Value = 34772.8
$ws.Range("N80").Value = -36768.8
$ws.Range("H83").Value = 34772.8
$ws.Range("J83").Value = 34772.8
$ws.Range("L83").Value = 104318.4
$ws.Range("N83").Value = -114302.4
$ws.Range("H97").Value = 1126.8
$ws.Range("I97").Value = 1035.3334
$ws.Range("J97").Value = 1950
$ws.Range("K97").Value = 1035.3334
$ws.Range("L97").Value = 1950
$ws.Range("M97").Value = -539.3334
$ws.Range("N97").Value = -2942

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 25251.334
$ws.Range("J35").Value = 25251.334
$ws.Range("L35").Value = 25251.334
$ws.Range("N35").Value = -25871.334
$ws.Range("H82").Value = 21575.77
$ws.Range("J82").Value = 33628.5
$ws.Range("L82").Value = 33628.5
$ws.Range("N82").Value = -34394.5
$ws.Range("H85").Value = 21575.77
$ws.Range("J85").Value = 33628.5
$ws.Range("L85").Value = 33628.5
$ws.Range("N85").Value = -36280.5
$ws.Range("H107").Value = 1131.5172
$ws.Range("I107").Value = 1031.8
$ws.Range("J107").Value = 1353.1111
$ws.Range("K107").Value = 1031.8
$ws.Range("L107").Value = 1353.1111
$ws.Range("M107").Value = 888.2
$ws.Range("N107").Value = -5193.1111

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1097.25
$ws.Range("I19").Value = 1097.25
$ws.Range("K19").Value = 1097.25
$ws.Range("M19").Value = -927.25
$ws.Range("H24").Value = 1097.25
$ws.Range("I24").Value = 1097.25
$ws.Range("K24").Value = 1097.25
$ws.Range("M24").Value = -927.25
$ws.Range("H132").Value = 13399.6
$ws.Range("I132").Value = 20000
$ws.Range("J132").Value = 8999.333000000001
$ws.Range("K132").Value = 60000
$ws.Range("L132").Value = 26997.999
$ws.Range("M132").Value = -57470
$ws.Range("N132").Value = -32057.999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 20649.75
$ws.Range("I12").Value = 11500
$ws.Range("K12").Value = 11500
$ws.Range("M12").Value = -11360
$ws.Range("H113").Value = 4045
$ws.Range("I113").Value = 4400
$ws.Range("J113").Value = 2980
$ws.Range("K113").Value = 4400
$ws.Range("L113").Value = 2980
$ws.Range("M113").Value = -2230
$ws.Range("N113").Value = -7320

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3187.4285
$ws.Range("I7").Value = 1524.5
$ws.Range("K7").Value = 1524.5
$ws.Range("M7").Value = -1412.5
$ws.Range("H122").Value = 9399.6
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 3187.4285
$ws.Range("I126").Value = 1524.5
$ws.Range("K126").Value = 4573.5
$ws.Range("M126").Value = -2103.5
$ws.Range("H132").Value = 6025.5625
$ws.Range("I132").Value = 4968
$ws.Range("J132").Value = 6660.1
$ws.Range("K132").Value = 14904
$ws.Range("L132").Value = 19980.3
$ws.Range("M132").Value = -12374
$ws.Range("N132").Value = -25040.3

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 39550
$ws.Range("J112").Value = 39550
$ws.Range("L112").Value = 39550
$ws.Range("N112").Value = -42504
$ws.Range("H113").Value = 464.08334
$ws.Range("I113").Value = 445.75
$ws.Range("K113").Value = 1337.25
$ws.Range("M113").Value = 832.75
$ws.Range("H115").Value = 37416.668
$ws.Range("J115").Value = 37416.668
$ws.Range("L115").Value = 37416.668
$ws.Range("N115").Value = -40550.668
$ws.Range("H122").Value = 3575.3333
$ws.Range("I122").Value = 2187.6667
$ws.Range("J122").Value = 4500.4443
$ws.Range("K122").Value = 6563.000100000001
$ws.Range("L122").Value = 13501.3329
$ws.Range("M122").Value = -4113.000100000001
$ws.Range("N122").Value = -18401.3329
$ws.Range("H123").Value = 33254.043
$ws.Range("J123").Value = 33254.043
$ws.Range("L123").Value = 33254.043
$ws.Range("N123").Value = -43054.043
$ws.Range("H125").Value = 38861
$ws.Range("J125").Value = 38861
$ws.Range("L125").Value = 38861
$ws.Range("N125").Value = -48701
$ws.Range("H126").Value = 446141.1
$ws.Range("I126").Value = 1463.2727
$ws.Range("J126").Value = 822406.9399999999
$ws.Range("K126").Value = 4389.8181
$ws.Range("L126").Value = 2467220.82
$ws.Range("M126").Value = -1919.8181
$ws.Range("N126").Value = -2472160.82
$ws.Range("H128").Value = 41657.5
$ws.Range("J128").Value = 41657.5
$ws.Range("L128").Value = 41657.5
$ws.Range("N128").Value = -51617.5
$ws.Range("H131").Value = 41584.445
$ws.Range("I131").Value = 40000
$ws.Range("J131").Value = 41782.5
$ws.Range("K131").Value = 40000
$ws.Range("L131").Value = 41782.5
$ws.Range("M131").Value = -34960
$ws.Range("N131").Value = -51862.5
$ws.Range("H132").Value = 3268.818
$ws.Range("I132").Value = 1422.7142
$ws.Range("K132").Value = 4268.142599999999
$ws.Range("M132").Value = -1738.142599999999
$ws.Range("H137").Value = 43134.285
$ws.Range("J137").Value = 43134.285
$ws.Range("L137").Value = 43134.285
$ws.Range("N137").Value = -53334.285
$ws.Range("H139").Value = 46721.11
$ws.Range("J139").Value = 46721.11
$ws.Range("L139").Value = 46721.11
$ws.Range("N139").Value = -57001.11
$ws.Range("H141").Value = 47763.43
$ws.Range("J141").Value = 47763.43
$ws.Range("L141").Value = 47763.43
$ws.Range("N141").Value = -58123.43
